$wb = $excel.ActiveWorkbook

# --- YDS sheet: append Week 16 play-by-play yardage logs ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = "10 1 0 4 4 6 6 3 12 0 1 12 0 3 3 1 4 9 6 6 4 3 3 12 3 3 3 3 3 4 -4 8 7 9 10 7 1 1 0 4 11 10 3 2 9 1 1 11 5 3 4 1 2 0 7 -3 1 -1 4 4 0 1 1 7 3 2 3 1 9 9 5 6 4 3 -1 5 3 9 0 2 24 1 5 2 3 -2 3 5 8 11 3 2 -2 2 9 9 5 5 3 8 3 8 1 0 1 9 7 1 13 9 12 0 0 7 4 7 -1 3 -2 0 4 4 11 1 2 12 1 4 3 8 10 5 -2 3 3 3 7 0 0 2 2 -1 5 7 1 2 6 2 5 1 1 4 5 0 4 0 1 12 1 6 4 6 5 2 0 6 0 0 4 6 0 7 19 3 5 11 4 -3 -3 4 1 -1 -3 3 9 6 8 7 2 9 2 10 0 11 0 0 4 2 14 6 6 -4 5 9 2 5 2 4 20 -1 -1 0 15 11 0 3 3 -2 5 8 3 20 5 12 -2 0 1 0 8 5 10 2 3 4 -1 2 5 3 4 7 5 9 4 6 7 4 5 0 1 4 3 1 3 15 6 42 4 5 2 18 3 6 2 4 8 -1 11 2 2 2 2 6 3 8 3 -2 5 4 3 1 7 2 9 2 13 4 1 2 10 3 -1 0 10 14 3 2 2 3 2 3 3 0 -1 0 3 12 2 2 5 3 -1 5 0 4 1 15 0 1 1 17 2 -1 4 7 26 -2 2 4 3 2 3 3 2 1 2 7 2 0 5 3 -1 2 6 1 4 1 2 4 5 11 6 40 1 5 6 3 1 -1 4 1 9 5 3 -1 0 3 5 5 2 5 1 7 1 3 14 9 3 -1 2 4 31 0 4 -1 19 3 13 -2 2 2 3 0 6 0 2 7 3 2 2 5 1 6 5 13 -2 2 3 4 4 4 9 4 2 4 0 2 7 1 0 3 -1 13 3 0 2 4 1 4 11 1 10 1 1 3 5 -2 5 4 8 16 6 13 5 6 4 5 5 6 19 5 28 3 23 3 1 -1 5 4 3 3 5 9 1 4 4 8 11 7 5 -1 6 5 7 6 4 9 5 8 1 2 10 3 0 0 13 2 -1 5 4 8 3 1 5 4 1 5 -6 2 2 19 0 8 7 5 7 1 4 0 1 2 6 1 11 5 5 47 14 0 14 0 0 10 5 2 18 8 1 5 2 15 4 4 6 6 17 4 5 3 6 6 3 5 12 0 0 14 0 4 3 7 -2 5 3 14 0 10 2 2 7 4 13 2 2 2 3 2 2 5 0 2 -1 9 4 1 12 3 4 -1 0 0 0 1 3 9 0 2 4 4 21 8 8 10 2 7 8 2 6 2 3 5 4 2 7 2 5 2 2 2 6 2 4 2 2 4 2 1 2 1 4 5 8 1 7 11 7 2 3 -2 7 2 5 10 9 3 5 4 -3 4 2 -3 3 2 1 9 9 2 2 2 1 0 3 1 3 15 22 12 7 2 3 -1 2 5 31 3 1 5 1 5 2 8 3 5 4 4 4 4 2 5 8 1 3 5 2 1 4 3 2 3 8 3 1 12 -5 1 1 0 3 2 1 2 33 2 1 3 4 -2 4 58 2 10 2 3 5 3 3 10 4 3 0 4 0 5 5 4 0 2 2 7 4 5 4 11 1 4 10 10 6 1 3 2 3 3 4 3 4 1 5 5 3 -1 7 8 4 13 17 2 3 2 0 2 3 10 7 9 12 2 5 0 0 2 1 -2 1 12 3 1 2 4 9 8 4 0 8 -2 0 13 0 11 4 9 3 -2 1 6 3 4 0 5 2 7 1 8 2"
$ws.Range("C2").Value = "9 10 24 4 5 33 6 19 3 18 9 12 12 14 9 11 12 11 11 14 2 5 7 4 15 8 10 10 9 11 9 -3 10 11 3 8 12 20 6 4 37 18 8 58 6 6 6 24 12 4 8 8 10 11 3 47 14 38 24 6 5 4 7 28 0 28 28 -6 13 22 40 5 10 -4 10 1 52 42 20 10 9 43 19 -2 10 8 3 6 2 10 6 2 4 8 8 14 12 25 18 43 10 13 9 19 5 6 17 0 7 3 7 6 7 10 11 14 13 18 7 16 11 12 9 26 12 5 14 9 16 5 12 10 12 6 11 16 6 18 10 14 3 13 4 10 8 5 26 6 21 11 22 11 8 8 9 -7 10 15 19 38 0 4 13 6 0 1 4 17 6 3 8 12 5 5 -2 26 -2 6 13 8 11 7 27 6 11 5 4 9 9 6 9 1 5 4 14 10 4 18 6 6 8 32 7 12 13 6 10 14 1 11 2 14 9 8 5 7 6 8 6 6 10 7 10 6 32 1 10 16 20 4 14 10 8 -4 16 8 5 14 9 14 2 32 20 6 21 6 4 14 19 4 13 3 6 5 10 32 10 6 9 7 -3 10 5 10 2 5 25 3 10 7 54 7 0 9 7 9 9 1 4 1 19 3 7 7 7 2 6 4 6 4 9 -4 0 2 14 13 11 13 12 13 6 6 15 9 4 6 13 13 9 11 6 4 19 19 12 8 12 -4 11 5 4 2 20 13 14 12 11 11 14 7 15 7 32 5 12 4 10 3 30 1 11 4 13 8 22 1 12 2 45 12 9 5 12 35 14 7 9 7 21 4 55 6 5 8 16 4 7 19 7 11 52 69 1 14 6 4 9 3 13 16 9 6 14 9 6 11 7 21 7 10 -3 5 20 7 10 8 3 16 3 10 5 3 5 7 28 5 16 11 8 12 7 8 7 22 2 5 22 6 8 3 4 13 -5 4 11 7 10 5 16 6 14 6 4 3 7 10 9 6 21 2 3 6 11 9 31 10 6 7 6 20 7 17 12 10 13 5 16 10 34 4 13 -4 1 11 7 7 5 12 4 8 8 44 7 8 4 12 19 9 8 6 9 6 15 15 13 11 7 5 22 18 12 10 14 18 4 7 20 35 6 7 6 23 8 6 8 6 2 49 6 8 -2 11 21 7 15 24 11 30 15 4 8 14 18 22 35 7 11 14 11 18 17 3 23 1 13 6 13 44 20 11 3 11 4 13 0 18 27 20 1 13 12 33 3 2 4 12 24 14 6 9 7 35 5 8 17 8 20 29 4 9 3 7 73 10 2 35 1 13 15 7 19 33 18 1 7 -3 15 5 5 0 32 17 21 4 18 5 4 24 9 32 5 10 7 4 20 11 4 37 11 13 5 9 12 6 6 18 5 7 4 9 8 16 9 14 9 19 7 22 11 6 8 -3 11 11 2 2 8 11 3 1 5 4 5 16 4 0 6 12 18 14 -2 3 19 5 19 7 11 16 1 6 15 17 10 10 9 4 51 5 13 5 7 7 15 7 41 12 7 6 5 4 32 32 11 5 3 17 10 2 2 4 9 5 1 41 7 1 21 10 10 5 17 1 4 20 0 12 2 10 1 3 4 13 25 9 4 3 4 18 7 15 3 12 24 7 4 0 12 13 9 14 20 12 5 11 5 10 2 4 9 11 4 2 19 4 9 9 6 5 11 3 3 6 7 15 9 3 7 3 10 16 10 1 11 7 9 7 7 7 5 10 24 6 5 7 14 22 4 10 9 5 40 8 9 25 16 13 1 18 11 12 8 13 8 61 9"
$ws.Range("B3").Value = "4 0 5 1 3 0 7 3 4 3 0 11 14 6 0 4 1 6 2 12 1 4 4 8 1 2 3 4 3 6 4 4 2 8 3 4 3 2 3 5 5 3 5 1 10 3 1 1 3 8 1 3 5 6 -2 7 4 2 1 3 1 6 9 4 -2 5 -2 4 6 4 -2 12 2 23 2 -1 1 4 11 1 1 0 0 9 1 4 -1 8 4 5 0 2 7 4 13 2 2 7 21 -2 9 2 3 23 12 8 4 10 12 28 10 2 4 6 12 2 6 8 9 3 14 4 1 5 11 -2 3 2 9 7 2 -5 6 50 8 4 9 2 9 3 2 3 8 2 4 0 1 3 2 5 1 2 4 1 7 -6 6 13 1 6 -2 4 4 6 3 15 3 0 10 4 3 1 -1 16 1 20 22 9 5 1 10 5 1 2 7 3 4 3 8 3 5 3 0 6 69 5 5 8 1 40 2 2 7 12 6 3 5 13 5 9 -1 4 5 8 5 12 5 3 1 5 3 10 2 1 3 3 3 4 -3 11 19 9 7 6 17 -4 1 7 4 11 3 2 4 2 2 4 0 7 -2 9 0 4 4 7 0 3 8 3 6 2 3 -2 1 0 6 1 4 6 2 3 4 4 4 2 11 3 5 3 1 1 -1 2 2 9 2 1 4 1 -1 5 13 5 5 8 7 6 16 1 6 0 4 2 2 6 1 5 5 3 3 3 6 13 6 10 1 2 -1 4 1 3 2 23 1 15 2 7 0 1 3 37 1 3 8 3 10 7 5 18 0 1 36 4 7 2 0 2 30 -1 4 2 2 7 7 5 5 14 21 9 -2 4 6 5 0 7 0 5 24 3 5 1 3 6 6 0 3 -2 2 9 2 4 0 9 0 3 5 9 4 3 4 11 3 0 1 5 4 0 5 1 8 2 5 1 8 4 9 17 -3 2 8 10 5 7 -1 0 1 1 9 0 5 10 1 6 3 3 1 0 3 1 4 12 -1 1 7 7 5 11 12 2 6 1 4 7 -1 2 2 7 20 1 -4 8 6 12 6 1 3 5 2 18 3 5 3 8 10 7 23 3 0 1 1 2 11 4 0 3 5 2 2 7 0 5 11 -1 1 4 5 8 5 3 5 3 3 5 2 5 1 7 0 6 6 1 20 -1 16 1 3 1 5 1 0 5 1 9 5 3 9 3 7 4 -2 24 2 3 3 1 3 7 9 8 2 4 14 12 2 -1 1 -1 14 7 4 11 -2 2 6 7 3 3 0 7 -2 7 7 4 6 5 1 4 3 4 1 14 2 4 1 0 0 -1 4 0 2 4 -2 7 3 4 21 4 4 5 3 1 1 1 11 -1 6 9 1 3 4 1 21 6 6 3 1 0 0 6 5 2 6 2 -1 15 0 2 3 2 0 0 13 19 4 3 -2 8 4 8 8 7 1 4 -4 -2 11 20 -3 7 13 1 -2 3 17 0 12 2 2 11 2 2 2 30 2 1 1 7 14 6 7 6 4 -4 2 0 1 10 3 2 0 2 1 -1 6 14 3 2 3 5 -3 4 1 -5 1 6 2 4 7 17 6 8 33 11 4 3 4 1 5 1 5 0 2 3 3 7 4 14 3 6 4 10 -2 -7 3 2 2 4 1 3 0 9 4 1 2 3 0 10 5 5 4 -4 6 6 3 22 2 4 4 4 2 4 5 21 -1 0 4 4 2 5 3 -1 3 2 9 9 3 -1 25 24 7 0 4 8 6 5 9 23 4 6 8 0 7 7 2 0 1 3 3 11 2 4 6 5 0 13 6 5 12 8 4 6 7 2 0 1 6 5 4 5 5 28 4 5 0 5 1 -3 0 6 6 3 5 1 3 2 0 6 3 15 15 4 1 0 8 3 10 3 0 4 1 0 4 4 2 3"
$ws.Range("C3").Value = "20 31 9 18 6 4 8 4 5 21 11 31 10 7 30 6 17 7 15 15 22 2 42 4 25 5 4 3 14 11 24 1 10 2 6 7 22 2 8 4 12 14 10 7 43 13 62 4 10 2 15 1 2 5 9 8 16 1 1 -1 8 11 16 6 6 4 21 5 29 11 37 16 5 16 10 4 11 1 11 9 13 12 6 7 14 11 1 2 5 12 17 18 7 9 3 11 21 16 12 16 -3 7 10 39 3 9 14 -1 2 1 5 6 13 6 6 11 80 60 7 1 9 5 26 52 6 7 15 13 5 22 12 8 6 8 32 6 2 18 6 3 6 4 6 4 9 2 15 9 2 -4 6 14 7 5 8 8 9 13 4 9 10 5 17 1 2 18 2 20 20 9 31 8 7 8 9 14 3 9 42 5 8 10 28 12 11 8 27 4 11 26 7 51 6 2 8 13 20 7 3 13 39 12 4 -1 7 28 4 6 4 6 9 8 5 5 12 5 9 15 9 20 8 10 8 1 3 4 12 8 4 38 4 2 20 9 2 10 11 16 6 7 21 11 5 7 9 10 0 6 10 5 4 9 3 1 7 2 5 4 16 9 5 5 14 12 13 5 22 19 9 -3 4 23 25 3 13 3 9 15 5 11 7 8 -1 23 3 26 49 3 5 81 19 5 17 23 10 15 43 10 8 9 22 17 6 10 14 9 8 14 16 9 15 21 10 10 4 2 2 18 33 7 16 3 9 21 13 10 10 19 28 16 8 5 27 5 19 17 2 10 47 20 2 2 11 10 1 5 7 9 7 9 7 3 11 8 4 6 20 24 22 11 7 0 5 14 -2 12 9 27 12 14 8 -1 6 10 7 8 42 8 -2 5 7 16 8 20 4 2 15 30 14 38 5 5 9 19 7 19 28 41 27 0 3 4 1 5 12 5 7 24 19 16 12 15 -4 9 9 29 9 3 13 6 5 17 1 8 5 9 39 9 15 2 11 4 7 13 8 6 5 5 8 38 5 9 18 28 6 35 6 2 26 5 13 13 16 14 13 13 11 7 13 9 14 22 20 27 0 8 7 3 5 11 3 11 75 14 31 3 20 3 5 7 -2 6 11 18 7 7 -4 -2 4 4 9 -3 6 13 6 3 -3 25 5 8 11 44 19 7 6 9 17 6 3 11 20 8 -3 16 40 17 18 17 14 0 8 13 21 9 1 10 8 12 0 8 37 13 23 7 11 8 24 14 22 10 12 6 6 4 8 8 8 12 56 5 12 12 21 17 13 31 3 11 -1 4 35 10 -5 16 54 9 9 4 30 7 8 6 12 -2 -6 25 27 9 24 9 1 9 30 6 1 5 13 13 70 14 5 12 2 3 19 9 43 -1 6 10 8 19 19 7 8 2 16 4 15 0 3 4 9 15 9 7 7 8 14 7 7 10 36 14 48 8 2 21 -4 32 15 -2 19 5 5 6 12"

# --- OFF sheet: updated season totals through Week 16 ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 398
$ws.Range("F2").Value = 117
$ws.Range("G2").Value = 107
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 71
$ws.Range("L2").Value = 573
$ws.Range("M2").Value = 386
$ws.Range("O2").Value = 50
$ws.Range("P2").Value = 26
$ws.Range("Q2").Value = 1011
$ws.Range("C3").Value = 385
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 240
$ws.Range("G3").Value = 90
$ws.Range("H3").Value = 39
$ws.Range("I3").Value = 116
$ws.Range("J3").Value = 109
$ws.Range("N3").Value = 30

# --- DEF sheet: updated season totals through Week 16 ---
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 402
$ws.Range("E2").Value = 25
$ws.Range("F2").Value = 138
$ws.Range("J2").Value = 71
$ws.Range("L2").Value = 518
$ws.Range("M2").Value = 311
$ws.Range("O2").Value = 41
$ws.Range("P2").Value = 21
$ws.Range("Q2").Value = 974
$ws.Range("C3").Value = 337
$ws.Range("F3").Value = 190
$ws.Range("G3").Value = 74
$ws.Range("H3").Value = 47
$ws.Range("I3").Value = 118
$ws.Range("J3").Value = 100
$ws.Range("N3").Value = 36

# --- ST sheet: updated totals + appended Week 16 return logs ---
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 182
$ws.Range("D2").Value = 106
$ws.Range("F2").Value = 350
$ws.Range("G2").Value = 336
$ws.Range("B3").Value = 138
$ws.Range("D3").Value = "35 32 35 43 48 42 52 41 40 48 39 44 41 45 32 47 44 45 49 37 44 36 40 39 50 40 52 58 52 47 48 51 52 41 53 42 52 42 52 39 31 56 48 39 48 53 49 54 38 47 62 44 56 49 40 57 51 40 36 63 61 51 50 54 36 41 55 55 46 38 51 60 52 46 47 34 48 53 60 55 53 44 58 56 59 60 59 59 54 31 37 37 43 31 59 35 50 50 32 36 50 42 47 54 62 51"
$ws.Range("B4").Value = "57 61 59 64 59 60 63 67 60 23 65 65 57 66 65 68 43 55 59 61 68 65 68 65 64 64 65 63 61 65 63 63 73 74 62 64 62 68 67 68 66 66 56 63"
$ws.Range("D4").Value = "0 0 0 0 0 0 13 0 0 0 0 0 -1 0 0 0 0 3 0 9 22 0 0 0 3 0 0 12 10 0 0 0 9 0 14 0 0 0 0 0 0 15 0 0 20 0 19 0 0 24 7 6 -4 -3 0 5 0 0 0 0 7 23 0 7 0 0 0 0 1 0 5 4 0 0 0 0 0 9 21 0 0 0 0 9 7 8 19 0 8 0 2 0 0 0 0 0 0 0 0 0 0 0 7 7 0 0"
$ws.Range("B5").Value = "20 17 20 25 9 17 21 32 20 8 25 24 14 25 20 29 20 18 23 28 24 23 23 19 15 24 22 18 23 20 41 31 33 21 23 24 18 38 27 39 28 20 15 0"
$ws.Range("D5").Value = "20 0 0 0 0 0 0 3 7 19 0 0 6 1 0 6 27 4 -5 3 0 1 0 0 0 0 0 4 0 3 83 6 18 10 20 9 0 0 1 -1 10 0 16 0 6 0 0 0 0 0 0 0 0 -2 0 3 12 0 5 8 0 0 0 0 0 0 0 0 17 0 -1 0 -1 0 0 7 0 -3 0 0 0 3 16 0 0 7 0 0 2 0 0 0 14 7 16 9 0 0 12 0 4 0 1 0 0 0 0 0 21 0 0 0 0 0"
$ws.Range("B6").Value = "21 27 13 22 12 28 31 7 26 26 28 13 15 28 67 30 15 20 14 23 64 16 15 66 29 22 17 60 30 20 25 33 21 18 22 22 24 23 24 14 25 22 27 22 54 35 29 21 18 18 0 11 22 21 15 21 29 20"

# --- TURNS sheet: updated totals through Week 16 ---
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 16
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 10
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 9

# --- PEN sheet: updated totals through Week 16 ---
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("D2").Value = 23
$ws.Range("D3").Value = 12

Write-Host "Applied Week 16 log + season sim updates"
